$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing order rows (2:43) down by one to make room for a new
# record at the top of the data, preserving each row's formatting (dates).
$ws.Range("A2:G43").Copy($ws.Range("A3:G44"))

# Fill in the newly opened row with the new order record.
$ws.Range("A2").Value = "4/29/2022"
$ws.Range("B2").Value = "Midwest"
$ws.Range("C2").Value = "Michael"
$ws.Range("D2").Value = "Paper"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1.29
$ws.Range("G2").Value = 3.87

# Match the saved selection from the edit.
$ws.Range("G3").Select()
